$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 32 and 33 swap places (WrappedliquidstakedEther2.0 <-> ImmutableX),
# along with their updated Price / Volume(1h) figures.
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.111"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.30%  "

$ws.Range("B33").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C33").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.883.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.13%  "

# Updated Price / Volume(1h) values for the remaining rows.
# Price (column D) cells are forced to text format before assignment so
# that values such as "316.50" or "0.3940" are not reinterpreted as
# numbers (which would drop trailing zeros / use scientific notation);
# the style is then reset to Normal so no stray cell formatting remains.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.592.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.695.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.92%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.06%  "

$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3940"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4025"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.525"
$ws.Range("D9").Style = "Normal"

$ws.Range("E10").Value = "  +0.13%  "

$ws.Range("E11").Value = "  +8.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08776"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.228"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001324"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.605"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.693.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "100.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07057"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.873"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.590.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.027"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.308"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.227"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.82%  "

$ws.Range("E31").Value = "  +15.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.311"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08526"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.955"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2736"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02783"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09071"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.465"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7705"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7191"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.555"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.221"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.93%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.356"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +12.71%  "

$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08022"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.26%  "
